# Updates cryptos list figures (price/volume columns) to match the
# latest scrape, including a rank swap between BitcoinCash and Uniswap
# (rows 21/22). For cells whose new text looks like a plain decimal
# number, force the cell to Text format first so Excel keeps the exact
# string instead of coercing it into a floating point number, then
# restore the default "Normal" style so no stray formatting is left
# behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.400.76'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '2.012.01'
$ws.Range('E3').Value = '  -0.61%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '259.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E6').Value = '  -1.70%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '57.27'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.386'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.24%  '
$ws.Range('E10').Value = '  -4.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.34'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.89%  '
$ws.Range('D13').Value = '2.309.54'
$ws.Range('E13').Value = '  -0.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.30'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.803'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.65%  '
$ws.Range('E16').Value = '  -5.19%  '
$ws.Range('D17').Value = '2.012.53'
$ws.Range('E17').Value = '  -0.76%  '
$ws.Range('D18').Value = '37.276.33'
$ws.Range('E18').Value = '  +0.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '70.01'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.95%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '232.78'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.63%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.15'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.60'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.42%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('E25').Value = '  -0.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.79'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.99'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.64'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.02%  '
$ws.Range('E29').Value = '  -6.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.34'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.97%  '
$ws.Range('E31').Value = '  -1.69%  '
$ws.Range('E32').Value = '  -2.78%  '
$ws.Range('E33').Value = '  -5.26%  '
$ws.Range('E34').Value = '  -0.80%  '
$ws.Range('E36').Value = '  +0.58%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('E38').Value = '  -3.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.41'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.52%  '
$ws.Range('E40').Value = '  +3.77%  '
$ws.Range('E41').Value = '  +0.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0213'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0930'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.30%  '
$ws.Range('D44').Value = '1.416.29'
$ws.Range('E44').Value = '  +1.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.85'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.74'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.71%  '
$ws.Range('E47').Value = '  -3.23%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.45%  '
$ws.Range('E49').Value = '  -6.22%  '
$ws.Range('D50').Value = '2.200.88'
$ws.Range('E50').Value = '  -0.67%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.96'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -9.99%  '
